# Updates res_bus/vm_pu.xlsx values for the "case with 380 kV" run (rows 2-25, B:F and I:N).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.057975094177015
$ws.Range("D2").Value = 1.057198732715922
$ws.Range("E2").Value = 1.063884259865603
$ws.Range("F2").Value = 1.074285903361631
$ws.Range("I2").Value = 1.047255178407985
$ws.Range("J2").Value = 1.062968722279131
$ws.Range("K2").Value = 1.059934358194967
$ws.Range("L2").Value = 1.066601692997784
$ws.Range("M2").Value = 1.076975518756298
$ws.Range("N2").Value = 1.064478259788983

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059255550931844
$ws.Range("D3").Value = 1.058195054297776
$ws.Range("E3").Value = 1.065098339828122
$ws.Range("F3").Value = 1.075749859333795
$ws.Range("I3").Value = 1.047642588113682
$ws.Range("J3").Value = 1.063900270894081
$ws.Range("K3").Value = 1.060744004898174
$ws.Range("L3").Value = 1.067629870247748
$ws.Range("M3").Value = 1.07825498305456
$ws.Range("N3").Value = 1.065411131309818

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060083402598871
$ws.Range("D4").Value = 1.058839047644678
$ws.Range("E4").Value = 1.065883590519318
$ws.Range("F4").Value = 1.076697225510287
$ws.Range("I4").Value = 1.04789169323551
$ws.Range("J4").Value = 1.064501843561879
$ws.Range("K4").Value = 1.06126659447615
$ws.Range("L4").Value = 1.068294259798602
$ws.Range("M4").Value = 1.079082448109097
$ws.Range("N4").Value = 1.066013558279805

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.060431269827482
$ws.Range("D5").Value = 1.059109618416007
$ws.Range("E5").Value = 1.066213631272102
$ws.Range("F5").Value = 1.077095523291868
$ws.Range("I5").Value = 1.047996041076183
$ws.Range("J5").Value = 1.064754459399422
$ws.Range("K5").Value = 1.061485980414837
$ws.Range("L5").Value = 1.068573353923771
$ws.Range("M5").Value = 1.07943021387318
$ws.Range("N5").Value = 1.066266532860812

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060489668888953
$ws.Range("D6").Value = 1.059155038876855
$ws.Range("E6").Value = 1.066269042037022
$ws.Range("F6").Value = 1.077162400820405
$ws.Range("I6").Value = 1.048013539522867
$ws.Range("J6").Value = 1.064796858047384
$ws.Range("K6").Value = 1.061522798114606
$ws.Range("L6").Value = 1.068620202508586
$ws.Range("M6").Value = 1.079488599475931
$ws.Range("N6").Value = 1.066308991719718

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.06008805144853
$ws.Range("D7").Value = 1.05884266366667
$ws.Range("E7").Value = 1.065888000845499
$ws.Range("F7").Value = 1.076702547481618
$ws.Range("I7").Value = 1.047893089012367
$ws.Range("J7").Value = 1.064505220145227
$ws.Range("K7").Value = 1.061269527140315
$ws.Range("L7").Value = 1.0682979899105
$ws.Range("M7").Value = 1.07908709536318
$ws.Range("N7").Value = 1.066016939658288

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058407975371499
$ws.Range("D8").Value = 1.057535588999023
$ws.Range("E8").Value = 1.06429463601546
$ws.Range("F8").Value = 1.074780638039057
$ws.Range("I8").Value = 1.047386432064586
$ws.Range("J8").Value = 1.063283793705032
$ws.Range("K8").Value = 1.060208253666867
$ws.Range("L8").Value = 1.066949360268923
$ws.Range("M8").Value = 1.077408011556388
$ws.Range("N8").Value = 1.064793778652449

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.055442031306236
$ws.Range("D9").Value = 1.055226958160917
$ws.Range("E9").Value = 1.061484191466523
$ws.Range("F9").Value = 1.071394497988453
$ws.Range("I9").Value = 1.04648152790952
$ws.Range("J9").Value = 1.061122176052087
$ws.Range("K9").Value = 1.058328064917145
$ws.Range("L9").Value = 1.064565817756664
$ws.Range("M9").Value = 1.07444577867602
$ws.Range("N9").Value = 1.062629091254495

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.053460862407867
$ws.Range("D10").Value = 1.05368411580283
$ws.Range("E10").Value = 1.059608545937636
$ws.Range("F10").Value = 1.069137174636916
$ws.Range("I10").Value = 1.045870047710484
$ws.Range("J10").Value = 1.059674689724243
$ws.Range("K10").Value = 1.057067698577123
$ws.Range("L10").Value = 1.062971861648673
$ws.Range("M10").Value = 1.072468410959027
$ws.Range("N10").Value = 1.061179549330037

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.052602026591676
$ws.Range("D11").Value = 1.053015128651346
$ws.Range("E11").Value = 1.05879585211024
$ws.Range("F11").Value = 1.068159699331536
$ws.Range("I11").Value = 1.045603306854969
$ws.Range("J11").Value = 1.059046358390034
$ws.Range("K11").Value = 1.056520279599021
$ws.Range("L11").Value = 1.062280455071779
$ws.Range("M11").Value = 1.071611537538323
$ws.Range("N11").Value = 1.060550325693264

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.05228286567567
$ws.Range("D12").Value = 1.052766495312108
$ws.Range("E12").Value = 1.058493898695817
$ws.Range("F12").Value = 1.067796610839076
$ws.Range("I12").Value = 1.045503930654731
$ws.Range("J12").Value = 1.058812730742398
$ws.Range("K12").Value = 1.056316690098449
$ws.Range("L12").Value = 1.062023450344897
$ws.Range("M12").Value = 1.071293153237459
$ws.Range("N12").Value = 1.06031636626757

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.052351333680812
$ws.Range("D13").Value = 1.052819834462126
$ws.Range("E13").Value = 1.058558672530046
$ws.Range("F13").Value = 1.067874495110857
$ws.Range("I13").Value = 1.045525260639297
$ws.Range("J13").Value = 1.058862855457851
$ws.Range("K13").Value = 1.056360372259211
$ws.Range("L13").Value = 1.062078587180274
$ws.Range("M13").Value = 1.071361452495344
$ws.Range("N13").Value = 1.060366562165868

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.052575647757364
$ws.Range("D14").Value = 1.052994579448013
$ws.Range("E14").Value = 1.058770894255455
$ws.Range("F14").Value = 1.068129686567961
$ws.Range("I14").Value = 1.045595098449488
$ws.Range("J14").Value = 1.059027051516432
$ws.Range("K14").Value = 1.056503456017754
$ws.Range("L14").Value = 1.062259214784912
$ws.Range("M14").Value = 1.071585221925447
$ws.Range("N14").Value = 1.060530991401686

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052713834917248
$ws.Range("D15").Value = 1.05310222674415
$ws.Range("E15").Value = 1.058901639993355
$ws.Range("F15").Value = 1.068286916895295
$ws.Range("I15").Value = 1.045638088461917
$ws.Range("J15").Value = 1.059128186577353
$ws.Range("K15").Value = 1.056591580948203
$ws.Range("L15").Value = 1.062370480729151
$ws.Range("M15").Value = 1.071723079835101
$ws.Range("N15").Value = 1.060632270085994

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.0535178395673
$ws.Range("D16").Value = 1.053728494543293
$ws.Range("E16").Value = 1.059662470343668
$ws.Range("F16").Value = 1.069202045220946
$ws.Range("I16").Value = 1.045887708882685
$ws.Range("J16").Value = 1.059716356877196
$ws.Range("K16").Value = 1.057103993529675
$ws.Range("L16").Value = 1.063017722187732
$ws.Range("M16").Value = 1.072525264560916
$ws.Range("N16").Value = 1.061221275655126

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05402190592523
$ws.Range("D17").Value = 1.054121086000596
$ws.Range("E17").Value = 1.060139575526244
$ws.Range("F17").Value = 1.069776067165525
$ws.Range("I17").Value = 1.046043761853824
$ws.Range("J17").Value = 1.060084880734986
$ws.Range("K17").Value = 1.057424966981766
$ws.Range("L17").Value = 1.063423392564581
$ws.Range("M17").Value = 1.073028274158885
$ws.Range("N17").Value = 1.061590322859066

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.054315825369406
$ws.Range("D18").Value = 1.054349988763137
$ws.Range("E18").Value = 1.060417812308958
$ws.Range("F18").Value = 1.070110880973796
$ws.Range("I18").Value = 1.046134595327902
$ws.Range("J18").Value = 1.060299684134388
$ws.Range("K18").Value = 1.057612024210016
$ws.Range("L18").Value = 1.063659896412851
$ws.Range("M18").Value = 1.073321607736377
$ws.Range("N18").Value = 1.061805431303932

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.054416028501724
$ws.Range("D19").Value = 1.054428023638305
$ws.Range("E19").Value = 1.060512675386734
$ws.Range("F19").Value = 1.070225043381883
$ws.Range("I19").Value = 1.046165535081072
$ws.Range("J19").Value = 1.060372901134181
$ws.Range("K19").Value = 1.05767577863154
$ws.Range("L19").Value = 1.063740518397695
$ws.Range("M19").Value = 1.073421616276102
$ws.Range("N19").Value = 1.061878752280263

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05396783412425
$ws.Range("D20").Value = 1.054078973897894
$ws.Range("E20").Value = 1.060088391916186
$ws.Range("F20").Value = 1.069714480439029
$ws.Range("I20").Value = 1.046027038467234
$ws.Range("J20").Value = 1.060045357185116
$ws.Range("K20").Value = 1.057390546256268
$ws.Range("L20").Value = 1.063379880080479
$ws.Range("M20").Value = 1.072974312579948
$ws.Range("N20").Value = 1.061550743181221

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.052509597057401
$ws.Range("D21").Value = 1.052943125338072
$ws.Range("E21").Value = 1.058708402559686
$ws.Range("F21").Value = 1.068054539376121
$ws.Range("I21").Value = 1.045574541158758
$ws.Range("J21").Value = 1.058978706447691
$ws.Range("K21").Value = 1.056461328441252
$ws.Range("L21").Value = 1.062206029605664
$ws.Range("M21").Value = 1.071519330292492
$ws.Range("N21").Value = 1.060482577677402

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.051591869425909
$ws.Range("D22").Value = 1.052228151009405
$ws.Range("E22").Value = 1.057840267978877
$ws.Range("F22").Value = 1.067010804623278
$ws.Range("I22").Value = 1.045288320416389
$ws.Range("J22").Value = 1.058306685690351
$ws.Range("K22").Value = 1.055875623130628
$ws.Range("L22").Value = 1.061466908108706
$ws.Range("M22").Value = 1.070603925088022
$ws.Range("N22").Value = 1.059809602573505

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.052078458702019
$ws.Range("D23").Value = 1.052607251135509
$ws.Range("E23").Value = 1.058300529289105
$ws.Range("F23").Value = 1.067564115804114
$ws.Range("I23").Value = 1.045440214742137
$ws.Range("J23").Value = 1.058663067961136
$ws.Range("K23").Value = 1.056186256701987
$ws.Range("L23").Value = 1.061858833408943
$ws.Range("M23").Value = 1.071089256988097
$ws.Range("N23").Value = 1.060166490947992

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.053992267133673
$ws.Range("D24").Value = 1.054098002819955
$ws.Range("E24").Value = 1.060111519742618
$ws.Range("F24").Value = 1.069742308840311
$ws.Range("I24").Value = 1.046034595631524
$ws.Range("J24").Value = 1.060063216640423
$ws.Range("K24").Value = 1.057406099998367
$ws.Range("L24").Value = 1.063399541859899
$ws.Range("M24").Value = 1.072998695690471
$ws.Range("N24").Value = 1.061568627999003

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.056209465494959
$ws.Range("D25").Value = 1.055824447181875
$ws.Range("E25").Value = 1.062211100891067
$ws.Range("F25").Value = 1.072269860532541
$ws.Range("I25").Value = 1.046716909466778
$ws.Range("J25").Value = 1.06168212445696
$ws.Range("K25").Value = 1.05881534710647
$ws.Range("L25").Value = 1.06518287709045
$ws.Range("M25").Value = 1.075212020081996
$ws.Range("N25").Value = 1.063189834850331
